# Apply updated TPM-derived values to the LR-pairs sheet (Col5a3-Sdc3)
# The sheet recomputes several NATMI metrics after the ligand/receptor
# average & total expression values were refreshed with new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    # Row 2
    $ws.Range("G2").Value = 0.6646009999999999
    $ws.Range("H2").Value = 1.993803
    $ws.Range("I2").Value = 0.006573892384547845
    $ws.Range("J2").Value = 0.006573892384547845
    $ws.Range("M2").Value = 28.19948866666667
    $ws.Range("N2").Value = 84.598466
    $ws.Range("O2").Value = 0.7357427920402423
    $ws.Range("P2").Value = 0.7357427920402422
    $ws.Range("Q2").Value = 18.74140836735533
    $ws.Range("R2").Value = 168.672675306198
    $ws.Range("S2").Value = 0.004836693937579318
    $ws.Range("T2").Value = 0.004836693937579317

    # Row 3
    $ws.Range("G3").Value = 0.6646009999999999
    $ws.Range("H3").Value = 1.993803
    $ws.Range("I3").Value = 0.006573892384547845
    $ws.Range("J3").Value = 0.006573892384547845
    $ws.Range("O3").Value = 0.2029336910395279
    $ws.Range("P3").Value = 0.2029336910395278
    $ws.Range("Q3").Value = 5.169283636092333
    $ws.Range("R3").Value = 46.523552724831
    $ws.Range("S3").Value = 0.001334064246092937
    $ws.Range("T3").Value = 0.001334064246092937

    # Row 4
    $ws.Range("G4").Value = 0.6646009999999999
    $ws.Range("H4").Value = 1.993803
    $ws.Range("I4").Value = 0.006573892384547845
    $ws.Range("J4").Value = 0.006573892384547845
    $ws.Range("M4").Value = 2.350402666666667
    $ws.Range("N4").Value = 7.051208000000001
    $ws.Range("O4").Value = 0.0613235169202299
    $ws.Range("P4").Value = 0.06132351692022989
    $ws.Range("Q4").Value = 1.562079962669333
    $ws.Range("R4").Value = 14.058719664024
    $ws.Range("S4").Value = 0.0004031342008755902
    $ws.Range("T4").Value = 0.0004031342008755901

    # Row 5
    $ws.Range("I5").Value = 0.8500764341604863
    $ws.Range("J5").Value = 0.8500764341604863
    $ws.Range("M5").Value = 28.19948866666667
    $ws.Range("N5").Value = 84.598466
    $ws.Range("O5").Value = 0.7357427920402423
    $ws.Range("P5").Value = 0.7357427920402422
    $ws.Range("Q5").Value = 2423.469789909362
    $ws.Range("R5").Value = 21811.22810918426
    $ws.Range("S5").Value = 0.6254376091168494
    $ws.Range("T5").Value = 0.6254376091168493

    # Row 6
    $ws.Range("I6").Value = 0.8500764341604863
    $ws.Range("J6").Value = 0.8500764341604863
    $ws.Range("O6").Value = 0.2029336910395279
    $ws.Range("P6").Value = 0.2029336910395278
    $ws.Range("S6").Value = 0.1725091484499077
    $ws.Range("T6").Value = 0.1725091484499076

    # Row 7
    $ws.Range("I7").Value = 0.8500764341604863
    $ws.Range("J7").Value = 0.8500764341604863
    $ws.Range("M7").Value = 2.350402666666667
    $ws.Range("N7").Value = 7.051208000000001
    $ws.Range("O7").Value = 0.0613235169202299
    $ws.Range("P7").Value = 0.06132351692022989
    $ws.Range("Q7").Value = 201.9940830885423
    $ws.Range("S7").Value = 0.05212967659372928
    $ws.Range("T7").Value = 0.05212967659372927

    # Row 8
    $ws.Range("I8").Value = 0.1433496734549659
    $ws.Range("J8").Value = 0.1433496734549659
    $ws.Range("M8").Value = 28.19948866666667
    $ws.Range("N8").Value = 84.598466
    $ws.Range("O8").Value = 0.7357427920402423
    $ws.Range("P8").Value = 0.7357427920402422
    $ws.Range("Q8").Value = 408.6733722416016
    $ws.Range("R8").Value = 3678.060350174414
    $ws.Range("S8").Value = 0.1054684889858136
    $ws.Range("T8").Value = 0.1054684889858136

    # Row 9
    $ws.Range("I9").Value = 0.1433496734549659
    $ws.Range("J9").Value = 0.1433496734549659
    $ws.Range("O9").Value = 0.2029336910395279
    $ws.Range("P9").Value = 0.2029336910395278
    $ws.Range("S9").Value = 0.02909047834352725
    $ws.Range("T9").Value = 0.02909047834352724

    # Row 10
    $ws.Range("I10").Value = 0.1433496734549659
    $ws.Range("J10").Value = 0.1433496734549659
    $ws.Range("M10").Value = 2.350402666666667
    $ws.Range("N10").Value = 7.051208000000001
    $ws.Range("O10").Value = 0.0613235169202299
    $ws.Range("P10").Value = 0.06132351692022989
    $ws.Range("Q10").Value = 34.06256741980356
    $ws.Range("R10").Value = 306.5631067782321
    $ws.Range("S10").Value = 0.00879070612562503
    $ws.Range("T10").Value = 0.008790706125625028

